# Apply the edit described by the commit:
#   "adding cells with multilines (multicells)"
#
# The underlying OOXML diff shows the product-name text in B3 being
# appended with extra text, and the sheet's active selection ending on
# that same cell (B3) instead of F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product name in B3 (shared string "Microscope Landless Land"
# -> "Microscope Landless Land fdghdsfyhsg").
$ws.Range("B3").Value = "Microscope Landless Land fdghdsfyhsg"

# Leave the selection on B3, matching the saved sheet view state.
$ws.Range("B3").Select() | Out-Null
